$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = "Update Page Layouts"
$ws.Range("E13").Value = "Update Account, Contact and Opportunity Page Layouts to keep related fields to the P2P Car Rental App"
$ws.Range("E14").Value = "Create a experience micro-site and pages for P2P Car Rental App"
$ws.Range("B19").Value = "Create Apex REST to send available cars to rent to external system"
$ws.Range("B20").Value = "Integrate Stripe for payments"
$ws.Range("B21").Value = "Show a map on home page to show available car locations"
$ws.Range("E24").Value = "Consider add additional relationship and you may need to  ensure correct vehicle record will be updated Hint: Trigger and Validation Rule"
$ws.Range("E25").Value = "Update InUse after creation and update (cover both scenario: Add vehicle or remove vehicle) Hint: Trigger"
$ws.Range("E26").Value = "Follow best practices for finding a solution. Coding or Point and click tools?"

$ws.Range("B24").Value = "Automatically update Last Odometer reading,  Odometer Reading date of Vehicle records based on last created Vehicle Renting  record"
$ws.Range("B28").Value = "If Vehicle Owner opt-in for Smart price, set Is smart price field true and price field will be unavailable automatically on Rental Listing"
$ws.Range("E28").Value = "Hint: Trigger for is Smart Price  and validation for price field"

$ws.Range("B29").Value = "Calculate renting amount automatically for Vehicle renting records based on the prices on rental listing"
$ws.Range("B30").Value = "Show surcharge fee automatically for Vehicle renting records if the Vehicle older than 15 years."
$ws.Range("B31").Value = "Calculate Snapp Car Fee for Vehicle renting records automatically based on Renting amount"
$ws.Range("B32").Value = "Show young driver fee for renters less than 25 years old for Vehicle renting records"
$ws.Range("B33").Value = "Calculate insurance fee 10% of the renting amount with a minimum of $ 5 automatically for Vehicle renting records"
$ws.Range("B34").Value = "Show One-off verification fee for contacts if their first rent on the platform for Vehicle renting records. It is $15"
$ws.Range("B35").Value = "Calculate extra mile fee if it is more than allowance on the listing for Vehicle renting records"
$ws.Range("B36").Value = "Calculate roadside assistance fee $2 daily for Vehicle renting records"
$ws.Range("B37").Value = "Calculate final amountautomatically for Vehicle renting records"
$ws.Range("E37").Value = "Tax included for all prices, apply long term discount if vehicle owner opt-in"
$ws.Range("B38").Value = "Validation rules for Vehicle"
$ws.Range("E38").Value = "VIN is required"
$ws.Range("B40").Value = "Throw an error for Vehicle renting if Vehicle is not active or rented "
$ws.Range("E40").Value = "Contact should be prevented to choose rental listing on Vehicle Renting object(Hint: Trigger)"
$ws.Range("B41").Value = "Send an email to vehicle owners once their rental listing is created and contacts once their vehicle renting created."
$ws.Range("B43").Value = "Rental Listings shpuld be inactive if end date is older than current date. "
$ws.Range("E43").Value = "Is active should be automatically false if listing end date is yesterday (at 12 am everyday) Hint: Schedule Apex"
$ws.Range("B44").Value = "Update Rental Listing Title automatically "
$ws.Range("E44").Value = "Rent + Vehicle Name + (If exist) Daily + (if exist) Hourly : `"Rent Tesla Model Y Daily Hourly`""
$ws.Range("B45").Value = "Deactivate all vehicles and rental listings automatically once Vehicle owner is not active"
$ws.Range("B46").Value = "Automatically delete all old Rental Listing records first Sundayof every month"
$ws.Range("E46").Value = "Automatically delete all records if it is older than 12 months. Notify admin user via email about the result  Hint: Batch Apex "
$ws.Range("B47").Value = "Calculate the total amount of vehicle renting for Vehicle owners automatically "
$ws.Range("B48").Value = "Update isKeyless field for Vehicle Renting records if a keyless assigned to the vehicle"
$ws.Range("B49").Value = "Client Interface"
$ws.Range("B50").Value = "Create a form to register as a car renter"
$ws.Range("B51").Value = "Create a list of available cars to rent"

# Numeric cell fills for rows 25-28 (Story Points / Priority columns)
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 3
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 3

# Update selection to match final view state
$ws.Range("B29").Select()
